{"js": "// This document's quiz table contains many short, distinct 'N+N=' / 'N-N=' cell\n// expressions plus a date title; the commit swaps the date and every expression for\n// a new value. Because every OLD string is unique in the document, we can safely\n// resolve all Range objects for the OLD strings first (single sync), then perform\n// all the text replacements against those already-resolved ranges (second sync).\n// Doing it this way (rather than re-searching text after each edit) avoids any\n// false-match hazard where a NEW value happens to equal an OLD value elsewhere in\n// the list (e.g. '15+43=' -> '28+15=' later collides with the original '28+15=' cell).\nconst pairs = [\n  [\"2023-11-03 Friday\", \"2023-11-04 Saturday\"],\n  [\"58+37=\", \"42-27=\"],\n  [\"1+3=\", \"77-37=\"],\n  [\"25+13=\", \"62+2=\"],\n  [\"49-36=\", \"68-42=\"],\n  [\"25+69=\", \"40-17=\"],\n  [\"97-82=\", \"88-62=\"],\n  [\"22+60=\", \"92-51=\"],\n  [\"43-8=\", \"71-7=\"],\n  [\"50-50=\", \"17+43=\"],\n  [\"67-55=\", \"67-28=\"],\n  [\"97-21=\", \"6+12=\"],\n  [\"16+41=\", \"46-39=\"],\n  [\"62-45=\", \"8-6=\"],\n  [\"61-40=\", \"71-53=\"],\n  [\"91-26=\", \"12+5=\"],\n  [\"27+50=\", \"9+36=\"],\n  [\"47-30=\", \"45-36=\"],\n  [\"12-11=\", \"30+68=\"],\n  [\"58-5=\", \"95-65=\"],\n  [\"22-4=\", \"72-66=\"],\n  [\"97+1=\", \"80-25=\"],\n  [\"74-45=\", \"47-9=\"],\n  [\"10+43=\", \"98-93=\"],\n  [\"57-28=\", \"11-1=\"],\n  [\"19+77=\", \"80-66=\"],\n  [\"42+52=\", \"63-58=\"],\n  [\"36+53=\", \"92-65=\"],\n  [\"61+7=\", \"17+18=\"],\n  [\"44+13=\", \"93-48=\"],\n  [\"94-55=\", \"63-13=\"],\n  [\"90+4=\", \"31+36=\"],\n  [\"11+87=\", \"81-64=\"],\n  [\"44-11=\", \"41+17=\"],\n  [\"51+19=\", \"15+51=\"],\n  [\"74+8=\", \"42+35=\"],\n  [\"82+2=\", \"70-36=\"],\n  [\"56+17=\", \"79-49=\"],\n  [\"60-26=\", \"26+33=\"],\n  [\"87-60=\", \"40+42=\"],\n  [\"58+33=\", \"78-19=\"],\n  [\"2+8=\", \"30+26=\"],\n  [\"73-51=\", \"20+20=\"],\n  [\"78-39=\", \"56-46=\"],\n  [\"85-80=\", \"13+5=\"],\n  [\"37-5=\", \"80-26=\"],\n  [\"63+19=\", \"52-19=\"],\n  [\"7+91=\", \"10+61=\"],\n  [\"17+39=\", \"81-65=\"],\n  [\"4+36=\", \"39-18=\"],\n  [\"31-6=\", \"18-9=\"],\n  [\"8+43=\", \"46-30=\"],\n  [\"60-16=\", \"98-51=\"],\n  [\"35+43=\", \"15+64=\"],\n  [\"29+63=\", \"3+63=\"],\n  [\"22+64=\", \"77-19=\"],\n  [\"43-10=\", \"58+41=\"],\n  [\"90-74=\", \"56-47=\"],\n  [\"59+6=\", \"31+63=\"],\n  [\"27+71=\", \"91-47=\"],\n  [\"14+59=\", \"29+41=\"],\n  [\"13+51=\", \"72-1=\"],\n  [\"2+85=\", \"96-21=\"],\n  [\"57-18=\", \"20+8=\"],\n  [\"63+29=\", \"99-27=\"],\n  [\"29-20=\", \"82+8=\"],\n  [\"98-59=\", \"61-18=\"],\n  [\"47+19=\", \"71+26=\"],\n  [\"65-31=\", \"66+29=\"],\n  [\"26+15=\", \"29+26=\"],\n  [\"1+43=\", \"59+3=\"],\n  [\"72-70=\", \"37+46=\"],\n  [\"48+34=\", \"84-52=\"],\n  [\"9-9=\", \"5+37=\"],\n  [\"70-64=\", \"19+58=\"],\n  [\"58-18=\", \"39-10=\"],\n  [\"41-16=\", \"22-9=\"],\n  [\"68-24=\", \"92-19=\"],\n  [\"44-6=\", \"93-62=\"],\n  [\"79-28=\", \"69-10=\"],\n  [\"53-22=\", \"44-25=\"],\n  [\"18+29=\", \"3+34=\"],\n  [\"93-57=\", \"39-27=\"],\n  [\"15+43=\", \"28+15=\"],\n  [\"19+31=\", \"86-65=\"],\n  [\"48-4=\", \"17-1=\"],\n  [\"94-23=\", \"82-36=\"],\n  [\"11+73=\", \"88-75=\"],\n  [\"16+40=\", \"81-24=\"],\n  [\"16+14=\", \"51-44=\"],\n  [\"78-57=\", \"30+31=\"],\n  [\"80-47=\", \"70-55=\"],\n  [\"67-27=\", \"73-63=\"],\n  [\"71+24=\", \"65-22=\"],\n  [\"76-66=\", \"25+33=\"],\n  [\"73-53=\", \"91-22=\"],\n  [\"85-6=\", \"44+30=\"],\n  [\"62-27=\", \"76-33=\"],\n  [\"58-12=\", \"72-33=\"],\n  [\"66-33=\", \"64-51=\"],\n  [\"28+15=\", \"63-9=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: locate every occurrence (search results are live Range objects, so\n// loading them now and editing them later is safe even once earlier edits have\n// changed surrounding text).\nconst searchResults = [];\nfor (const [oldText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\n// Phase 2: replace the text of each previously-found range with its new value.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const results = searchResults[i];\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# This document's quiz table contains many short, distinct 'N+N=' / 'N-N=' cell\n# expressions plus a date title; the commit swaps the date and every expression for\n# a new value. All OLD strings are unique in the document, but some NEW strings equal\n# an OLD string that appears LATER in the list (e.g. '15+43=' -> '28+15=' collides with\n# the original '28+15=' cell that comes after it). A plain 'search the whole document\n# each time' find/replace would therefore sometimes hit an already-edited cell instead\n# of the intended original one. To avoid that, we walk the document once, left to right:\n# each Find is scoped to start right where the previous replacement ended, so it can\n# only ever match the next not-yet-processed occurrence.\n$pairs = @(\n  ,@('2023-11-03 Friday', '2023-11-04 Saturday')\n  ,@('58+37=', '42-27=')\n  ,@('1+3=', '77-37=')\n  ,@('25+13=', '62+2=')\n  ,@('49-36=', '68-42=')\n  ,@('25+69=', '40-17=')\n  ,@('97-82=', '88-62=')\n  ,@('22+60=', '92-51=')\n  ,@('43-8=', '71-7=')\n  ,@('50-50=', '17+43=')\n  ,@('67-55=', '67-28=')\n  ,@('97-21=', '6+12=')\n  ,@('16+41=', '46-39=')\n  ,@('62-45=', '8-6=')\n  ,@('61-40=', '71-53=')\n  ,@('91-26=', '12+5=')\n  ,@('27+50=', '9+36=')\n  ,@('47-30=', '45-36=')\n  ,@('12-11=', '30+68=')\n  ,@('58-5=', '95-65=')\n  ,@('22-4=', '72-66=')\n  ,@('97+1=', '80-25=')\n  ,@('74-45=', '47-9=')\n  ,@('10+43=', '98-93=')\n  ,@('57-28=', '11-1=')\n  ,@('19+77=', '80-66=')\n  ,@('42+52=', '63-58=')\n  ,@('36+53=', '92-65=')\n  ,@('61+7=', '17+18=')\n  ,@('44+13=', '93-48=')\n  ,@('94-55=', '63-13=')\n  ,@('90+4=', '31+36=')\n  ,@('11+87=', '81-64=')\n  ,@('44-11=', '41+17=')\n  ,@('51+19=', '15+51=')\n  ,@('74+8=', '42+35=')\n  ,@('82+2=', '70-36=')\n  ,@('56+17=', '79-49=')\n  ,@('60-26=', '26+33=')\n  ,@('87-60=', '40+42=')\n  ,@('58+33=', '78-19=')\n  ,@('2+8=', '30+26=')\n  ,@('73-51=', '20+20=')\n  ,@('78-39=', '56-46=')\n  ,@('85-80=', '13+5=')\n  ,@('37-5=', '80-26=')\n  ,@('63+19=', '52-19=')\n  ,@('7+91=', '10+61=')\n  ,@('17+39=', '81-65=')\n  ,@('4+36=', '39-18=')\n  ,@('31-6=', '18-9=')\n  ,@('8+43=', '46-30=')\n  ,@('60-16=', '98-51=')\n  ,@('35+43=', '15+64=')\n  ,@('29+63=', '3+63=')\n  ,@('22+64=', '77-19=')\n  ,@('43-10=', '58+41=')\n  ,@('90-74=', '56-47=')\n  ,@('59+6=', '31+63=')\n  ,@('27+71=', '91-47=')\n  ,@('14+59=', '29+41=')\n  ,@('13+51=', '72-1=')\n  ,@('2+85=', '96-21=')\n  ,@('57-18=', '20+8=')\n  ,@('63+29=', '99-27=')\n  ,@('29-20=', '82+8=')\n  ,@('98-59=', '61-18=')\n  ,@('47+19=', '71+26=')\n  ,@('65-31=', '66+29=')\n  ,@('26+15=', '29+26=')\n  ,@('1+43=', '59+3=')\n  ,@('72-70=', '37+46=')\n  ,@('48+34=', '84-52=')\n  ,@('9-9=', '5+37=')\n  ,@('70-64=', '19+58=')\n  ,@('58-18=', '39-10=')\n  ,@('41-16=', '22-9=')\n  ,@('68-24=', '92-19=')\n  ,@('44-6=', '93-62=')\n  ,@('79-28=', '69-10=')\n  ,@('53-22=', '44-25=')\n  ,@('18+29=', '3+34=')\n  ,@('93-57=', '39-27=')\n  ,@('15+43=', '28+15=')\n  ,@('19+31=', '86-65=')\n  ,@('48-4=', '17-1=')\n  ,@('94-23=', '82-36=')\n  ,@('11+73=', '88-75=')\n  ,@('16+40=', '81-24=')\n  ,@('16+14=', '51-44=')\n  ,@('78-57=', '30+31=')\n  ,@('80-47=', '70-55=')\n  ,@('67-27=', '73-63=')\n  ,@('71+24=', '65-22=')\n  ,@('76-66=', '25+33=')\n  ,@('73-53=', '91-22=')\n  ,@('85-6=', '44+30=')\n  ,@('62-27=', '76-33=')\n  ,@('58-12=', '72-33=')\n  ,@('66-33=', '64-51=')\n  ,@('28+15=', '63-9=')\n)\n\n$d = $word.ActiveDocument\n$searchStart = 0\nforeach ($p in $pairs) {\n  $oldText = $p[0]\n  $newText = $p[1]\n  $rng = $d.Range($searchStart, $d.Content.End)\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n  if (-not $found) { throw \"Not found: $oldText\" }\n  $rng.Text = $newText\n  $searchStart = $rng.End\n}\n"}
